# Insert a new data row at row 304 (pushing the existing rows 304:376 down
# to 305:377) and populate the new row with the new market-price record.
# This matches the diff, which shows the whole table from row 304 onward
# shifted down by one row, with a brand-new record occupying row 304 and
# the sheet's dimension growing from A1:R376 to A1:R377.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 304:376 down to 305:377, leaving row 304 empty.
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new record.
$ws.Range("A304").Value = 4
$ws.Range("B304").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C304").Value = 'Los Lagos'
$ws.Range("D304").Value = 45015
$ws.Range("E304").Value = 10
$ws.Range("F304").Value = 100112032
$ws.Range("G304").Value = 'Zapallo italiano'
$ws.Range("H304").Value = 'Sin especificar'
$ws.Range("I304").Value = 'Primera'
$ws.Range("J304").Value = 120
$ws.Range("K304").Value = 11000
$ws.Range("L304").Value = 12000
$ws.Range("M304").Value = 11500
$ws.Range("N304").Value = '$/caja 50 unidades'
$ws.Range("O304").Value = 'Región Metropolitana'
$ws.Range("P304").Value = 230
$ws.Range("Q304").Value = 50
$ws.Range("R304").Value = 'Hortaliza'

# Make sure the date cell keeps the date number format used by the rest of
# column D (Insert() normally carries this down from the row above, but set
# it explicitly to be safe).
$ws.Range("D304").NumberFormat = "YYYY-MM-DD HH:MM:SS"
